$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1032.4
$ws.Range("I28").Value = 355.6
$ws.Range("J28").Value = 2386
$ws.Range("K28").Value = 355.6
$ws.Range("L28").Value = 2386
$ws.Range("M28").Value = 129.4
$ws.Range("N28").Value = -3356
$ws.Range("H43").Value = 3535.8
$ws.Range("I43").Value = 672.75
$ws.Range("J43").Value = 5444.5
$ws.Range("K43").Value = 672.75
$ws.Range("L43").Value = 5444.5
$ws.Range("M43").Value = -603.75
$ws.Range("N43").Value = -5582.5
$ws.Range("H132").Value = 2097.44
$ws.Range("I132").Value = 1758.0869
$ws.Range("K132").Value = 5274.2607
$ws.Range("M132").Value = -2744.2607
$ws.Range("H137").Value = 20002102
$ws.Range("I137").Value = 32259872
$ws.Range("J137").Value = 2580.8948
$ws.Range("K137").Value = 96779616
$ws.Range("L137").Value = 7742.6844
$ws.Range("M137").Value = -96777066
$ws.Range("N137").Value = -12842.6844
$ws.Range("H139").Value = 44496.09
$ws.Range("J139").Value = 44496.09
$ws.Range("L139").Value = 44496.09
$ws.Range("N139").Value = -54776.09

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6061.722
$ws.Range("I61").Value = 5423.6665
$ws.Range("K61").Value = 5423.6665
$ws.Range("M61").Value = -5211.6665
$ws.Range("H74").Value = 3048.75
$ws.Range("I74").Value = 2127.158
$ws.Range("K74").Value = 2127.158
$ws.Range("M74").Value = -1253.158
$ws.Range("H77").Value = 3048.75
$ws.Range("I77").Value = 2127.158
$ws.Range("K77").Value = 10635.79
$ws.Range("M77").Value = -6267.789999999999
$ws.Range("H102").Value = 688.125
$ws.Range("I102").Value = 700.6
$ws.Range("J102").Value = 667.3333
$ws.Range("K102").Value = 700.6
$ws.Range("L102").Value = 667.3333
$ws.Range("M102").Value = 921.4
$ws.Range("N102").Value = -3911.3333
$ws.Range("H122").Value = 1959.7
$ws.Range("I122").Value = 1612
$ws.Range("J122").Value = 2655.1
$ws.Range("K122").Value = 4836
$ws.Range("L122").Value = 7965.299999999999
$ws.Range("M122").Value = -2386
$ws.Range("N122").Value = -12865.3
$ws.Range("H130").Value = 59999
$ws.Range("J130").Value = 59999
$ws.Range("L130").Value = 59999
$ws.Range("N130").Value = -70039
$ws.Range("H136").Value = 6061.722
$ws.Range("I136").Value = 5423.6665
$ws.Range("K136").Value = 16270.9995
$ws.Range("M136").Value = -13720.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2019.1
$ws.Range("I94").Value = 3517.6
$ws.Range("J94").Value = 520.6
$ws.Range("K94").Value = 3517.6
$ws.Range("L94").Value = 520.6
$ws.Range("M94").Value = -3066.6
$ws.Range("N94").Value = -1422.6
$ws.Range("H99").Value = 2077.2222
$ws.Range("I99").Value = 1836.875
$ws.Range("K99").Value = 1836.875
$ws.Range("M99").Value = -338.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2702.1765
$ws.Range("I16").Value = 2026.8
$ws.Range("K16").Value = 2026.8
$ws.Range("M16").Value = -1739.8
$ws.Range("H31").Value = 35381.53
$ws.Range("I31").Value = 1740.5
$ws.Range("J31").Value = 78634.28999999999
$ws.Range("K31").Value = 1740.5
$ws.Range("L31").Value = 78634.28999999999
$ws.Range("M31").Value = -1445.5
$ws.Range("N31").Value = -79224.28999999999
$ws.Range("H34").Value = 35381.53
$ws.Range("I34").Value = 1740.5
$ws.Range("J34").Value = 78634.28999999999
$ws.Range("K34").Value = 1740.5
$ws.Range("L34").Value = 78634.28999999999
$ws.Range("M34").Value = -1538.5
$ws.Range("N34").Value = -79038.28999999999
$ws.Range("H58").Value = 3191.889
$ws.Range("I58").Value = 1843
$ws.Range("K58").Value = 1843
$ws.Range("M58").Value = -1640
$ws.Range("H113").Value = 2702.1765
$ws.Range("I113").Value = 2026.8
$ws.Range("K113").Value = 2026.8
$ws.Range("M113").Value = 143.2
$ws.Range("H134").Value = 3888
$ws.Range("I134").Value = 2336
$ws.Range("K134").Value = 7008
$ws.Range("M134").Value = -4473
$ws.Range("H136").Value = 3191.889
$ws.Range("I136").Value = 1843
$ws.Range("K136").Value = 5529
$ws.Range("M136").Value = -2979

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 13891066
$ws.Range("I68").Value = 25001160
$ws.Range("K68").Value = 75003480
$ws.Range("M68").Value = -75002669
$ws.Range("H71").Value = 13891066
$ws.Range("I71").Value = 25001160
$ws.Range("K71").Value = 225010440
$ws.Range("M71").Value = -225006384
$ws.Range("H112").Value = 125007380
$ws.Range("J112").Value = 4500
$ws.Range("L112").Value = 13500
$ws.Range("N112").Value = -15716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2291.4285
$ws.Range("I113").Value = 1604.2307
$ws.Range("K113").Value = 1604.2307
$ws.Range("M113").Value = 565.7692999999999
$ws.Range("H122").Value = 3035.2
$ws.Range("I122").Value = 2478.6956
$ws.Range("J122").Value = 4863.7144
$ws.Range("K122").Value = 7436.0868
$ws.Range("L122").Value = 14591.1432
$ws.Range("M122").Value = -4986.0868
$ws.Range("N122").Value = -19491.1432
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2512904.2
$ws.Range("I40").Value = 3490004.8
$ws.Range("K40").Value = 3490004.8
$ws.Range("M40").Value = -3489868.8
$ws.Range("H44").Value = 99215.336
$ws.Range("J44").Value = 99437.5
$ws.Range("L44").Value = 99437.5
$ws.Range("N44").Value = -100349.5
$ws.Range("H93").Value = 1848.0646
$ws.Range("I93").Value = 1994.0834
$ws.Range("K93").Value = 1994.0834
$ws.Range("M93").Value = -746.0834
$ws.Range("H115").Value = 78249.5
$ws.Range("J115").Value = 78249.5
$ws.Range("L115").Value = 78249.5
$ws.Range("N115").Value = -80599.5
$ws.Range("H122").Value = 91023.06
$ws.Range("I122").Value = 111766.22
$ws.Range("J122").Value = 5745.6665
$ws.Range("K122").Value = 335298.66
$ws.Range("L122").Value = 17236.9995
$ws.Range("M122").Value = -332848.66
$ws.Range("N122").Value = -22136.9995
$ws.Range("H132").Value = 3607.6592
$ws.Range("I132").Value = 3191.0645
$ws.Range("K132").Value = 9573.193499999999
$ws.Range("M132").Value = -7043.193499999999
$ws.Range("H136").Value = 4277.558
$ws.Range("I136").Value = 3028.682
$ws.Range("K136").Value = 9086.045999999998
$ws.Range("M136").Value = -6536.045999999998
$ws.Range("H139").Value = 62932.332
$ws.Range("J139").Value = 63923.875
$ws.Range("L139").Value = 63923.875
$ws.Range("N139").Value = -74203.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 921.4545000000001
$ws.Range("I100").Value = 923.1111
$ws.Range("J100").Value = 914
$ws.Range("K100").Value = 1846.2222
$ws.Range("L100").Value = 1828
$ws.Range("M100").Value = -1305.2222
$ws.Range("N100").Value = -2910
$ws.Range("H113").Value = 372.15384
$ws.Range("I113").Value = 236.625
$ws.Range("K113").Value = 709.875
$ws.Range("M113").Value = 1460.125
$ws.Range("H126").Value = 1907.9722
$ws.Range("I126").Value = 1385.7587
$ws.Range("K126").Value = 4157.2761
$ws.Range("M126").Value = -1687.2761
$ws.Range("H132").Value = 3223.8928
$ws.Range("I132").Value = 2530.76
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 7592.280000000001
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -5062.280000000001
$ws.Range("N132").Value = -32060
$ws.Range("H139").Value = 66566
$ws.Range("J139").Value = 68879.2
$ws.Range("L139").Value = 68879.2
$ws.Range("N139").Value = -79159.2
